$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Weekly Quantity" ---
$ws1 = $wb.Worksheets.Item("Weekly Quantity")

# Update B18 value from 20 to 4
$ws1.Range("B18").Value = 4

# Delete rows 19 and 20 entirely (shrinks used range to A1:B18)
$ws1.Rows.Item(19).Resize(2).Delete()

# --- Sheet 2: "Monthly Trend" ---
$ws2 = $wb.Worksheets.Item("Monthly Trend")

# Update B6 value from 230 to 41
$ws2.Range("B6").Value = 41
